# Append a new bonus record (row 28) to the "Worksheet" sheet, matching the
# formatting already used by the existing data rows (date format on column B,
# the "Ordem de Serviço" number style on column C).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate the last existing row (27) into a new row 28 so the new row
# inherits the same cell styles (date format on B, OS-number style on C)
# already present on row 27, instead of creating brand-new style entries.
$ws.Rows("27:27").Copy()
$ws.Rows("28:28").Insert(-4121)   # xlShiftDown
$excel.CutCopyMode = 0

# Overwrite the copied values with the new record's data.
$ws.Range("A28").Value = 26
$ws.Range("B28").Value = 45835        # 27/06/2025
$ws.Range("C28").Value = 69514093     # Ordem de Serviço
$ws.Range("D28").Value = 5            # Bonificação (R$)
$ws.Range("E28").Value = "Bruno"      # Técnico
$ws.Range("F28").Value = "Bruno fez sozinho"

# Match the author's final selection/active cell.
$ws.Range("D28").Select()
